# Data/g19.10.xlsx — "programada a coleta de dados VDE e organização das
# tabelas derivadas da fonte": refresh the source table with the newly
# collected values (Brasil / Nordeste / Sergipe), swap the Ano/Variável
# column order, and drop the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: B and C swap (Ano now before Variável) -------------------
$ws.Cells.Item(1, 1).Value = "Região"
$ws.Cells.Item(1, 2).Value = "Ano"
$ws.Cells.Item(1, 3).Value = "Variável"
$ws.Cells.Item(1, 4).Value = "Valor"
$ws.Cells.Item(1, 5).Value = "Posição relativamente às demais UF"

# --- Data rows --------------------------------------------------------------
$data = @(
    @("Brasil",   "01/01/2021", "Furto de veículo", 72.35437189743179, $null),
    @("Brasil",   "01/01/2022", "Furto de veículo", 81.75347312601198, $null),
    @("Brasil",   "01/01/2023", "Furto de veículo", 75.22949751208755, $null),
    @("Brasil",   "01/01/2024", "Furto de veículo", 70.4151726425975,  $null),
    @("Nordeste", "01/01/2021", "Furto de veículo", 39.17205746864656, $null),
    @("Nordeste", "01/01/2022", "Furto de veículo", 52.62425387813053, $null),
    @("Nordeste", "01/01/2023", "Furto de veículo", 54.12414104503122, $null),
    @("Nordeste", "01/01/2024", "Furto de veículo", 50.86320200298211, $null),
    @("Sergipe",  "01/01/2021", "Furto de veículo", 28.47750533632532, 26),
    @("Sergipe",  "01/01/2022", "Furto de veículo", 38.50982254945205, 26),
    @("Sergipe",  "01/01/2023", "Furto de veículo", 36.06223913262109, 24),
    @("Sergipe",  "01/01/2024", "Furto de veículo", 32.19900404429514, 26)
)

$lastRow = 1 + $data.Length

# The "Ano" column holds dd/mm/yyyy-looking text, not real dates — format the
# column as Text first so Excel doesn't auto-convert the literals to serials.
$anoRange = $ws.Range("B2:B$lastRow")
$anoRange.NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    if ($null -eq $row[4]) {
        # Rows 8-9 are already a blank column-E cell in the source sheet —
        # leave them untouched rather than round-tripping through
        # ClearContents(), which would needlessly re-flag them as changed.
        $current = $ws.Cells.Item($rowIndex, 5).Value2
        if ($null -ne $current -and $current -ne "") {
            $ws.Cells.Item($rowIndex, 5).ClearContents()
        }
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    }
    $rowIndex++
}

# Drop the temporary Text format now that the literal strings are committed.
$anoRange.Style = "Normal"

# --- Drop the now-unused trailing rows (old data had 19 rows, new has 13) --
$ws.Rows("14:19").Delete()

Write-Host "g19.10 table refreshed"
